$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 6bef14aa row refreshed
$wsOverview.Range("G2").Value = "2016-08-29 00:47:16"

# zh-cn sheet: Correspond Handoff/Handback datetimes refreshed for 6bef14aa row
$wsZhCn.Range("H2").Value = "2016-08-29 00:47:12"
$wsZhCn.Range("K2").Value = "2016-08-29 00:47:30"

# de-de sheet: Correspond Handoff/Handback datetimes refreshed for 6bef14aa row
$wsDeDe.Range("H2").Value = "2016-08-29 00:47:16"
$wsDeDe.Range("K2").Value = "2016-08-29 00:47:36"
